# Auto-generated edit script: updates market price / profit columns (H-N)
# across multiple leve-profit worksheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 239.04347
$ws.Range("I19").Value2 = 225
$ws.Range("J19").Value2 = 249.84616
$ws.Range("K19").Value2 = 225
$ws.Range("L19").Value2 = 249.84616
$ws.Range("M19").Value2 = -50
$ws.Range("N19").Value2 = -599.8461600000001
$ws.Range("H114").Value2 = 40712
$ws.Range("J114").Value2 = 40712
$ws.Range("L114").Value2 = 40712
$ws.Range("N114").Value2 = -49390
$ws.Range("H132").Value2 = 14008.311
$ws.Range("I132").Value2 = 2068.9219
$ws.Range("J132").Value2 = 90420.39999999999
$ws.Range("K132").Value2 = 6206.7657
$ws.Range("L132").Value2 = 271261.2
$ws.Range("M132").Value2 = -3676.7657
$ws.Range("N132").Value2 = -276321.2
$ws.Range("H135").Value2 = 15152940
$ws.Range("I135").Value2 = 1245
$ws.Range("J135").Value2 = 33334974
$ws.Range("K135").Value2 = 11205
$ws.Range("L135").Value2 = 300014766
$ws.Range("M135").Value2 = -8670
$ws.Range("N135").Value2 = -300019836
$ws.Range("H137").Value2 = 4257.0464
$ws.Range("I137").Value2 = 1383.3334
$ws.Range("J137").Value2 = 4723.054
$ws.Range("K137").Value2 = 4150.0002
$ws.Range("L137").Value2 = 14169.162
$ws.Range("M137").Value2 = -1600.0002
$ws.Range("N137").Value2 = -19269.162
$ws.Range("H138").Value2 = 1466.57
$ws.Range("I138").Value2 = 752.2
$ws.Range("J138").Value2 = 1851.2307
$ws.Range("K138").Value2 = 2256.6
$ws.Range("L138").Value2 = 5553.6921
$ws.Range("M138").Value2 = 2883.4
$ws.Range("N138").Value2 = -15833.6921
$ws.Range("H141").Value2 = 6046.933
$ws.Range("I141").Value2 = 3586.25
$ws.Range("J141").Value2 = 8859.143
$ws.Range("K141").Value2 = 10758.75
$ws.Range("L141").Value2 = 26577.429
$ws.Range("M141").Value2 = -5578.75
$ws.Range("N141").Value2 = -36937.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2001.7333
$ws.Range("I2").Value2 = 2240
$ws.Range("K2").Value2 = 2240
$ws.Range("M2").Value2 = -2127
$ws.Range("H32").Value2 = 28508.857
$ws.Range("I32").Value2 = 32095.312
$ws.Range("J32").Value2 = 17032.2
$ws.Range("K32").Value2 = 32095.312
$ws.Range("L32").Value2 = 17032.2
$ws.Range("M32").Value2 = -31808.312
$ws.Range("N32").Value2 = -17606.2
$ws.Range("H45").Value2 = 1322.2858
$ws.Range("I45").Value2 = 1135.3334
$ws.Range("J45").Value2 = 1462.5
$ws.Range("K45").Value2 = 1135.3334
$ws.Range("L45").Value2 = 1462.5
$ws.Range("M45").Value2 = -758.3334
$ws.Range("N45").Value2 = -2216.5
$ws.Range("H97").Value2 = 1350.05
$ws.Range("I97").Value2 = 1275.5555
$ws.Range("J97").Value2 = 1411
$ws.Range("K97").Value2 = 1275.5555
$ws.Range("L97").Value2 = 1411
$ws.Range("M97").Value2 = -779.5554999999999
$ws.Range("N97").Value2 = -2403
$ws.Range("H102").Value2 = 15648.066
$ws.Range("I102").Value2 = 1960.7273
$ws.Range("K102").Value2 = 1960.7273
$ws.Range("M102").Value2 = -338.7273
$ws.Range("H110").Value2 = 1448.6154
$ws.Range("I110").Value2 = 1448.6154
$ws.Range("K110").Value2 = 1448.6154
$ws.Range("M110").Value2 = 596.3846000000001
$ws.Range("H116").Value2 = 2001.7333
$ws.Range("I116").Value2 = 2240
$ws.Range("K116").Value2 = 2240
$ws.Range("M116").Value2 = 54
$ws.Range("H122").Value2 = 2067.6765
$ws.Range("I122").Value2 = 2046.4642
$ws.Range("K122").Value2 = 6139.392599999999
$ws.Range("M122").Value2 = -3689.392599999999
$ws.Range("H132").Value2 = 10418092
$ws.Range("I132").Value2 = 15625975
$ws.Range("J132").Value2 = 2324.75
$ws.Range("K132").Value2 = 46877925
$ws.Range("L132").Value2 = 6974.25
$ws.Range("M132").Value2 = -46875395
$ws.Range("N132").Value2 = -12034.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2001.7333
$ws.Range("I3").Value2 = 2240
$ws.Range("K3").Value2 = 2240
$ws.Range("M3").Value2 = -2126
$ws.Range("H105").Value2 = 2426.2942
$ws.Range("I105").Value2 = 1969.6666
$ws.Range("J105").Value2 = 2675.3635
$ws.Range("K105").Value2 = 1969.6666
$ws.Range("L105").Value2 = 2675.3635
$ws.Range("M105").Value2 = -222.6666
$ws.Range("N105").Value2 = -6169.363499999999
$ws.Range("H107").Value2 = 2074.2083
$ws.Range("I107").Value2 = 1975.0769
$ws.Range("J107").Value2 = 2191.3635
$ws.Range("K107").Value2 = 1975.0769
$ws.Range("L107").Value2 = 2191.3635
$ws.Range("M107").Value2 = -55.07690000000002
$ws.Range("N107").Value2 = -6031.363499999999
$ws.Range("H134").Value2 = 2840.25
$ws.Range("I134").Value2 = 1035.8096
$ws.Range("J134").Value2 = 4146.9136
$ws.Range("K134").Value2 = 3107.4288
$ws.Range("L134").Value2 = 12440.7408
$ws.Range("M134").Value2 = -572.4288000000001
$ws.Range("N134").Value2 = -17510.7408

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 151465.88
$ws.Range("I31").Value2 = 1786.0416
$ws.Range("J31").Value2 = 202061.88
$ws.Range("K31").Value2 = 1786.0416
$ws.Range("L31").Value2 = 202061.88
$ws.Range("M31").Value2 = -1491.0416
$ws.Range("N31").Value2 = -202651.88
$ws.Range("H34").Value2 = 151465.88
$ws.Range("I34").Value2 = 1786.0416
$ws.Range("J34").Value2 = 202061.88
$ws.Range("K34").Value2 = 1786.0416
$ws.Range("L34").Value2 = 202061.88
$ws.Range("M34").Value2 = -1584.0416
$ws.Range("N34").Value2 = -202465.88
$ws.Range("H132").Value2 = 51018.31
$ws.Range("I132").Value2 = 1654.909
$ws.Range("K132").Value2 = 4964.727000000001
$ws.Range("M132").Value2 = -2434.727000000001
$ws.Range("H134").Value2 = 484623.8
$ws.Range("I134").Value2 = 984.0952
$ws.Range("J134").Value2 = 1754178.1
$ws.Range("K134").Value2 = 2952.2856
$ws.Range("L134").Value2 = 5262534.300000001
$ws.Range("M134").Value2 = -417.2856000000002
$ws.Range("N134").Value2 = -5267604.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 2021.1875
$ws.Range("I97").Value2 = 1488.5333
$ws.Range("J97").Value2 = 10011
$ws.Range("K97").Value2 = 1488.5333
$ws.Range("L97").Value2 = 10011
$ws.Range("M97").Value2 = -992.5333000000001
$ws.Range("N97").Value2 = -11003
$ws.Range("H122").Value2 = 1189.6666
$ws.Range("I122").Value2 = 1300
$ws.Range("J122").Value2 = 1134.5
$ws.Range("K122").Value2 = 3900
$ws.Range("L122").Value2 = 3403.5
$ws.Range("M122").Value2 = -1450
$ws.Range("N122").Value2 = -8303.5
$ws.Range("H126").Value2 = 5310.643
$ws.Range("I126").Value2 = 8025.25
$ws.Range("J126").Value2 = 1691.1666
$ws.Range("K126").Value2 = 24075.75
$ws.Range("L126").Value2 = 5073.4998
$ws.Range("M126").Value2 = -21605.75
$ws.Range("N126").Value2 = -10013.4998
$ws.Range("H132").Value2 = 4332.5
$ws.Range("I132").Value2 = 1412.8462
$ws.Range("J132").Value2 = 7252.154
$ws.Range("K132").Value2 = 4238.5386
$ws.Range("L132").Value2 = 21756.462
$ws.Range("M132").Value2 = -1708.5386
$ws.Range("N132").Value2 = -26816.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2736.9062
$ws.Range("I7").Value2 = 2521.6365
$ws.Range("J7").Value2 = 3210.5
$ws.Range("K7").Value2 = 2521.6365
$ws.Range("L7").Value2 = 3210.5
$ws.Range("M7").Value2 = -2409.6365
$ws.Range("N7").Value2 = -3434.5
$ws.Range("H40").Value2 = 2166.8
$ws.Range("I40").Value2 = 2166.8333
$ws.Range("J40").Value2 = 2166.6667
$ws.Range("K40").Value2 = 2166.8333
$ws.Range("L40").Value2 = 2166.6667
$ws.Range("M40").Value2 = -2030.8333
$ws.Range("N40").Value2 = -2438.6667
$ws.Range("H61").Value2 = 4714.8335
$ws.Range("I61").Value2 = 4657.8
$ws.Range("K61").Value2 = 4657.8
$ws.Range("M61").Value2 = -4455.8
$ws.Range("H68").Value2 = 3445.3635
$ws.Range("I68").Value2 = 3271.2856
$ws.Range("J68").Value2 = 3750
$ws.Range("K68").Value2 = 3271.2856
$ws.Range("L68").Value2 = 3750
$ws.Range("M68").Value2 = -2522.2856
$ws.Range("N68").Value2 = -5248
$ws.Range("H71").Value2 = 3445.3635
$ws.Range("I71").Value2 = 3271.2856
$ws.Range("J71").Value2 = 3750
$ws.Range("K71").Value2 = 16356.428
$ws.Range("L71").Value2 = 18750
$ws.Range("M71").Value2 = -12612.428
$ws.Range("N71").Value2 = -26238
$ws.Range("H82").Value2 = 5556580
$ws.Range("I82").Value2 = 1142.2727
$ws.Range("J82").Value2 = 20834034
$ws.Range("K82").Value2 = 1142.2727
$ws.Range("L82").Value2 = 20834034
$ws.Range("M82").Value2 = -781.2727
$ws.Range("N82").Value2 = -20834756
$ws.Range("H85").Value2 = 5556580
$ws.Range("I85").Value2 = 1142.2727
$ws.Range("J85").Value2 = 20834034
$ws.Range("K85").Value2 = 1142.2727
$ws.Range("L85").Value2 = 20834034
$ws.Range("M85").Value2 = 105.7273
$ws.Range("N85").Value2 = -20836530
$ws.Range("H113").Value2 = 4714.8335
$ws.Range("I113").Value2 = 4657.8
$ws.Range("K113").Value2 = 4657.8
$ws.Range("M113").Value2 = -2487.8
$ws.Range("H126").Value2 = 2736.9062
$ws.Range("I126").Value2 = 2521.6365
$ws.Range("J126").Value2 = 3210.5
$ws.Range("K126").Value2 = 7564.9095
$ws.Range("L126").Value2 = 9631.5
$ws.Range("M126").Value2 = -5094.9095
$ws.Range("N126").Value2 = -14571.5
$ws.Range("H132").Value2 = 3672.1035
$ws.Range("I132").Value2 = 1928.7142
$ws.Range("J132").Value2 = 5299.2666
$ws.Range("K132").Value2 = 5786.142599999999
$ws.Range("L132").Value2 = 15897.7998
$ws.Range("M132").Value2 = -3256.142599999999
$ws.Range("N132").Value2 = -20957.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 2103
$ws.Range("I132").Value2 = 1655.7778
$ws.Range("J132").Value2 = 2606.125
$ws.Range("K132").Value2 = 4967.3334
$ws.Range("L132").Value2 = 7818.375
$ws.Range("M132").Value2 = -2437.3334
$ws.Range("N132").Value2 = -12878.375
$ws.Range("H136").Value2 = 19920.629
$ws.Range("I136").Value2 = 36459.25
$ws.Range("J136").Value2 = 2109.8076
$ws.Range("K136").Value2 = 109377.75
$ws.Range("L136").Value2 = 6329.4228
$ws.Range("M136").Value2 = -106827.75
$ws.Range("N136").Value2 = -11429.4228

Write-Host "Updated 250 cells across 7 sheets."
